# Adicionando comentarios no Scraping e Comecando analise dos dados:
# plotando os valores dos modelos agrupados por marca
#
# Appends the newly scraped Tabela Fipe rows (474-481) to the bottom of
# the data sheet, matching the columns:
# A Valor | B Marca | C Modelo | D AnoModelo | E Combustivel |
# F CodigoFipe | G MesReferencia | H Autenticacao | I TipoVeiculo |
# J SiglaCombustivel | K DataConsulta

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 17:58"),
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:00"),
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:00"),
    @("R$ 13.629,00", "GM - Chevrolet", "Celta 1.0/Super/N.Piq.1.0 MPFi VHC 8V 3p", 2001, "Gasolina", "004202-1", "fevereiro de 2001 ", "hg3826wm6jc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:03"),
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:05"),
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:05"),
    @("R$ 13.629,00", "GM - Chevrolet", "Celta 1.0/Super/N.Piq.1.0 MPFi VHC 8V 3p", 2001, "Gasolina", "004202-1", "fevereiro de 2001 ", "hg3826wm6jc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:07"),
    @("R$ 29.620,00", "GM - Chevrolet", "Astra 2.0/ CD/ GLS 2.0 MPFI 16V 3p", 2000, "Gasolina", "004169-6", "fevereiro de 2001 ", "pjmjmk4kqmc", 1, "G", "sexta-feira, 2 de fevereiro de 2024 18:08")
)

$startRow = $ws.UsedRange.Rows.Count + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
